$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109 (shifting existing rows 109+ down by one).
$ws.Rows("109:109").Insert()

# The freshly inserted row picks up "default" blank styles; copy the
# formatting from what is now row 111 (the row that used to be 110,
# i.e. the alternating "odd" row style) down onto the new row 109 so the
# new row matches the expected alternating banded style.
$ws.Range("A111:D111").Copy()
$ws.Range("A109:D109").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row with the new door entry.
$ws.Range("A109").Value = "H210"
$ws.Range("B109").Value = "H210"
$ws.Range("C109").Value = "Αίθουσα Διδασκαλίας"
$ws.Range("D109").Value = "Lecture Room - H4"

# Update the selected cell to reflect where the editor ended up.
$ws.Range("D110").Select() | Out-Null
